$d = $word.ActiveDocument

# wdFindContinue = 1, wdReplaceAll = 2
$wdFindContinue = 1
$wdReplaceAll = 2

function Replace-AllText($findText, $replaceText) {
    $d.Content.Find.Execute($findText, $true, $true, $false, $false, $false, `
                             $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll) | Out-Null
}

# [Muziki] -> [Music]
Replace-AllText "[Muziki]" "[Music]"

# kwa mfano -> for example (3 occurrences in the document)
Replace-AllText "kwa mfano" "for example"

# [SItisha] -> [PAUSE] (2 occurrences in the document)
Replace-AllText "[SItisha]" "[PAUSE]"
